$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("To Do")

# New ToDo item: "Turn Chart.component into a directive" - Open
$ws.Range("B20").Value = "Turn Chart.component into a directive"
$ws.Range("C20").Value = "Open"

# Existing row 28 ("Editing of Sheets by central user") gets a Status of "Open"
$ws.Range("C28").Value = "Open"

# New ToDo item: "Inject js lib" - Open
$ws.Range("B30").Value = "Inject js lib"
$ws.Range("C30").Value = "Open"

# Update the view: scroll so row 7 is the top-left visible row, and select C28
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C28").Select()
